$d = $word.ActiveDocument

# --- 1) Paragraph 1: merge "Replace "/"MongoDb"/" with "/"postgresql" runs into a
#        single run, and highlight the whole paragraph (incl. the paragraph mark) green.
$p1xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:rPr>
                <w:highlight w:val="green"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:highlight w:val="green"/>
              </w:rPr>
              <w:t>Replace MongoDb with postgresql</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
[void]$d.Paragraphs(1).Range.InsertXML($p1xml)

# --- 2) Paragraph 2: keep the text, but highlight the whole paragraph (incl. the
#        paragraph mark) green as well.
$p2xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:rPr>
                <w:highlight w:val="green"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:highlight w:val="green"/>
              </w:rPr>
              <w:t>Replace stripe with Khalti (Future)</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
[void]$d.Paragraphs(2).Range.InsertXML($p2xml)

# --- 3) Insert a new bullet right after paragraph 2 (same list formatting, no
#        highlight) for the new "In Report" item.
# InsertParagraphAfter() would carry over paragraph 2's green highlight onto the
# new paragraph mark, so stamp the freshly-created paragraph's XML directly
# (same technique as above) to get plain, unhighlighted formatting.
[void]$d.Paragraphs(2).Range.InsertParagraphAfter()
$p3xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>In Report: Add steps to add new product and model.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
[void]$d.Paragraphs(3).Range.InsertXML($p3xml)
